# Tugas kelompok: nambahin NPM (student ID numbers) after each member's name
# on Slide 1, "Content Placeholder 2" shape.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Paragraph 6: "M. Ken Fahriza" -> "M. Ken Fahriza - 54415678"
$para6 = $tr.Paragraphs(6, 1)
$null = $para6.InsertAfter(" - 54415678")

# Paragraph 7: "Silvano Satria" -> "Silvano Satria - 56415570"
$para7 = $tr.Paragraphs(7, 1)
$null = $para7.InsertAfter(" ")
$null = $para7.InsertAfter("- 56415570")

# Paragraph 8: "Yudha Patria" -> "Yudha " + "Patria - 57415312"
$para8 = $tr.Paragraphs(8, 1)
$part2 = $para8.Characters(7, 6)
$part2.Text = "Patria - 57415312"
